$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja3")

# Insert a new column before column C (pushes old C -> D, old D -> E)
$ws.Columns("C").Insert()

# Header
$ws.Range("C4").Value = "TIPO DATO O VALOR QUE RETORNA CADA FUNCION"

# Fill new column C per row
$ws.Range("C6").Value = "BOOLEAN"
$ws.Range("C7").Value = "BOOLEAN"
$ws.Range("C8").Value = "BOOLEAN"
$ws.Range("C10").Value = "STRING CCLCLAVE"
$ws.Range("C11").Value = "BOOLEAN"
$ws.Range("C12").Value = "BOOLEAN"
$ws.Range("C13").Value = "BOOLEAN"
$ws.Range("C14").Value = "BOOLEAN"
$ws.Range("C15").Value = "BOOLEAN"
$ws.Range("C16").Value = "BOOLEAN"
$ws.Range("C17").Value = "STRING ""S"" / ""N"""
$ws.Range("C18").Value = "STRING TEXTO DEL CAMPO"
$ws.Range("C19").Value = "STRING ""S"" / ""N"""
$ws.Range("C20").Value = "STRING ""PREPAGADO"" / ""POR COBRAR"""
$ws.Range("C21").Value = "STRING CCLCLAVE Y DIECLAVE"
$ws.Range("C22").Value = "BOOLEAN"
$ws.Range("C30").Value = "STRING ID_TIPO_TARIFA"
$ws.Range("C31").Value = "STRING CANTIDAD"

# Column widths per diff (C,D,E)
$ws.Columns("C").ColumnWidth = 47.140625
$ws.Columns("D").ColumnWidth = 90.28515625
$ws.Columns("E").ColumnWidth = 10.5703125
